# Applies the edits described by the commit:
#   "Moved all scripts into "scripts" folder and started documentation
#    for "Temperature responses.R""
#
# Four textual changes inside word/document.xml:
#  1. "... in the GitHub repo." -> '... in the "Model Parameters" folder.'
#  2. 'User-defined location for climate data (from "Habitat temperature
#      parameters.csv")' -> 'User-defined location for climate data or all = TRUE'
#  3. Append a new explanatory sentence right after "...nonlinear regression)."
#  4. The "Read in, and then find selected population in, " text is left
#     textually unchanged but gets re-merged into a single run (matches the
#     diff's run simplification).

$d = $word.ActiveDocument

$lq = [char]0x201C   # “  (left double quotation mark)
$rq = [char]0x201D   # ”  (right double quotation mark)

# --- Edit 1: GitHub repo -> "Model Parameters" folder ----------------------
$find1 = "the GitHub repo."
$repl1 = "the " + $lq + "Model Parameters" + $rq + " folder."
$ok1 = $d.Content.Find.Execute($find1, $true, $false, $false, $false, $false, `
                                $true, 1, $false, $repl1, 2)
Write-Output "Edit 1 (GitHub repo -> Model Parameters folder): $ok1"

# --- Edit 2: drop the "(from ...)" parenthetical, add "or all = TRUE" ------
$find2 = "User-defined location for climate data (from " + $lq + `
         "Habitat temperature parameters.csv" + $rq + ")"
$repl2 = "User-defined location for climate data or all = TRUE"
$ok2 = $d.Content.Find.Execute($find2, $true, $false, $false, $false, $false, `
                                $true, 1, $false, $repl2, 2)
Write-Output "Edit 2 (User-defined location ...): $ok2"

# --- Edit 3: append guidance about all = TRUE / all = FALSE ----------------
$find3 = "nonlinear regression)."
$repl3 = "nonlinear regression). Set all = TRUE if the script is to be run " + `
         "for all locations in " + $lq + "Habitat temperature parameters.csv" + $rq + `
         " or set all = FALSE if the script is to be run just for the specified location."
$ok3 = $d.Content.Find.Execute($find3, $true, $false, $false, $false, $false, `
                                $true, 1, $false, $repl3, 2)
Write-Output "Edit 3 (append Set all = TRUE/FALSE sentence): $ok3"

# --- Edit 4: normalize the "Read in, ..." run split (text unchanged) -------
$find4 = "Read in, and then find selected population in, "
$ok4 = $d.Content.Find.Execute($find4, $true, $false, $false, $false, $false, `
                                $true, 1, $false, $find4, 2)
Write-Output "Edit 4 (Read in, and then find selected population in, ): $ok4"
